$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.175.02"
$ws.Range("E2").Value = "  +2.55%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.086.94"
$ws.Range("E3").Value = "  +3.32%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "250.51"
$ws.Range("E5").Value = "  +2.89%  "
$ws.Range("E6").Value = "  +1.29%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "54.89"
$ws.Range("E8").Value = "  +25.39%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "61.68"
$ws.Range("E9").Value = "  +2.21%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.375"
$ws.Range("E10").Value = "  +5.68%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0745"
$ws.Range("E11").Value = "  +5.02%  "
$ws.Range("E12").Value = "  +8.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.04"
$ws.Range("E13").Value = "  +6.71%  "
$ws.Range("E14").Value = "  +3.33%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.830"
$ws.Range("E15").Value = "  +4.48%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.094.64"
$ws.Range("E16").Value = "  +4.06%  "
$ws.Range("E17").Value = "  +6.96%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.206.93"
$ws.Range("E18").Value = "  +2.68%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "72.73"
$ws.Range("E19").Value = "  +3.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.60"
$ws.Range("E20").Value = "  +16.63%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0847"
$ws.Range("E21").Value = "  +5.52%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "240.37"
$ws.Range("E22").Value = "  +2.66%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.20"
$ws.Range("E23").Value = "  +7.57%  "
$ws.Range("E24").Value = "  +0.21%  "
$ws.Range("E25").Value = "  +1.48%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "172.15"
$ws.Range("E26").Value = "  +2.68%  "
$ws.Range("E27").Value = "  +7.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.67"
$ws.Range("E28").Value = "  +5.54%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.00"
$ws.Range("E29").Value = "  +4.81%  "
$ws.Range("E30").Value = "  +2.82%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.72"
$ws.Range("E31").Value = "  +6.50%  "
$ws.Range("E32").Value = "  +29.99%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.52"
$ws.Range("E33").Value = "  +5.56%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0618"
$ws.Range("E34").Value = "  +7.95%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0898"
$ws.Range("E35").Value = "  +3.38%  "
$ws.Range("E36").Value = "  -0.15%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.15"
$ws.Range("E37").Value = "  +5.60%  "
$ws.Range("E38").Value = "  -3.10%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.21"
$ws.Range("E39").Value = "  +4.80%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.35"
$ws.Range("E40").Value = "  +3.51%  "
$ws.Range("B41").Value = "FTXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.72"
$ws.Range("E41").Value = "  +145.30%  "
$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "18.21"
$ws.Range("E42").Value = "  +19.37%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0226"
$ws.Range("E43").Value = "  +7.15%  "
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.16"
$ws.Range("E44").Value = "  +6.33%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "98.42"
$ws.Range("E45").Value = "  +3.80%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0946"
$ws.Range("E46").Value = "  +16.30%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.79"
$ws.Range("E47").Value = "  +0.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.322.65"
$ws.Range("E48").Value = "  +1.75%  "
$ws.Range("E49").Value = "  +5.56%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.34"
$ws.Range("E50").Value = "  +7.85%  "
$ws.Range("B51").Value = "FraxShare"
$ws.Range("C51").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.97"
$ws.Range("E51").Value = "  +15.35%  "
